$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Row 7: V7 goes 1 -> 0, and the running total X7 goes 13 -> 12
# ------------------------------------------------------------------
$ws.Range("V7").Value = 0
$ws.Range("X7").Value = 12

# ------------------------------------------------------------------
# 2) Column X (the rolling-weight total) is recomputed for every
#    policy row from 23 to 221 because the denominator used in the
#    underlying calculation changed (more weeks now known). The new
#    values come in contiguous same-value runs, so apply them as
#    row-range fills.
# ------------------------------------------------------------------
$xRuns = @(
    @{Start=23;  End=23;  Val=0.1111111111083333},
    @{Start=24;  End=31;  Val=0.1944444444416667},
    @{Start=32;  End=34;  Val=0.5833333333249999},
    @{Start=35;  End=37;  Val=0.6111111110999999},
    @{Start=38;  End=97;  Val=0.8333333333166668},
    @{Start=98;  End=98;  Val=0.7222222222083333},
    @{Start=99;  End=110; Val=0.6944444444333332},
    @{Start=111; End=112; Val=0.6111111110999999},
    @{Start=113; End=221; Val=0.3333333333333333}
)

foreach ($run in $xRuns) {
    for ($r = $run.Start; $r -le $run.End; $r++) {
        $ws.Cells.Item($r, 24).Value = $run.Val
    }
}

# ------------------------------------------------------------------
# 3) Twelve new policy/date rows are appended (222-233), one per new
#    calendar day from 9/30/2020 through 10/11/2020. Each new row
#    repeats the same weight pattern already seen on row 221.
# ------------------------------------------------------------------
$newDates = @(
    "9/30/2020", "10/1/2020", "10/2/2020", "10/3/2020",
    "10/4/2020", "10/5/2020", "10/6/2020", "10/7/2020",
    "10/8/2020", "10/9/2020", "10/10/2020", "10/11/2020"
)

$rowPattern = @(0,0,1,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,1)  # B..W

$firstNewRow = 222
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $firstNewRow + $i

    # Write the date label with a leading quote so Excel keeps it as
    # literal text instead of auto-converting it to a date serial.
    $ws.Cells.Item($r, 1).Value = "'" + $newDates[$i]

    for ($c = 0; $c -lt $rowPattern.Length; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $rowPattern[$c]
    }

    $ws.Cells.Item($r, 24).Value = 0.3333333333333333
}

# Re-apply column A's label formatting (bold / bordered / centered,
# matching every other date-label cell in the column) to the newly
# added cells - assigning the text value above resets it.
$ws.Range("A221").Copy()
$ws.Range("A222:A233").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
